# Apply updated TPM-derived NATMI ligand-receptor edge statistics to Sheet1.
# Final data spans all 5 sending clusters x 4 target clusters (20 rows total).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bgn"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = [double]"3.0"
$ws.Range("F2").Value = [double]"1.0"
$ws.Range("G2").Value = [double]"194.0250906666667"
$ws.Range("H2").Value = [double]"582.075272"
$ws.Range("I2").Value = [double]"0.04827281913040843"
$ws.Range("J2").Value = [double]"0.04827281913040844"
$ws.Range("K2").Value = [double]"3.0"
$ws.Range("L2").Value = [double]"1.0"
$ws.Range("M2").Value = [double]"4.618552666666667"
$ws.Range("N2").Value = [double]"13.855658"
$ws.Range("O2").Value = [double]"0.78434648953826"
$ws.Range("P2").Value = [double]"0.78434648953826"
$ws.Range("Q2").Value = [double]"896.1150998987752"
$ws.Range("R2").Value = [double]"8065.035899088976"
$ws.Range("S2").Value = [double]"0.03786261622505122"
$ws.Range("T2").Value = [double]"0.03786261622505122"

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bgn"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = [double]"3.0"
$ws.Range("F3").Value = [double]"1.0"
$ws.Range("G3").Value = [double]"194.0250906666667"
$ws.Range("H3").Value = [double]"582.075272"
$ws.Range("I3").Value = [double]"0.04827281913040843"
$ws.Range("J3").Value = [double]"0.04827281913040844"
$ws.Range("K3").Value = [double]"2.0"
$ws.Range("L3").Value = [double]"0.6666666666666666"
$ws.Range("M3").Value = [double]"0.6792986666666666"
$ws.Range("N3").Value = [double]"2.037896"
$ws.Range("O3").Value = [double]"0.1153620112191035"
$ws.Range("P3").Value = [double]"0.1153620112191036"
$ws.Range("Q3").Value = [double]"131.8009853897458"
$ws.Range("R3").Value = [double]"1186.208868507712"
$ws.Range("S3").Value = [double]"0.005568849502099933"
$ws.Range("T3").Value = [double]"0.005568849502099934"

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bgn"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = [double]"3.0"
$ws.Range("F4").Value = [double]"1.0"
$ws.Range("G4").Value = [double]"194.0250906666667"
$ws.Range("H4").Value = [double]"582.075272"
$ws.Range("I4").Value = [double]"0.04827281913040843"
$ws.Range("J4").Value = [double]"0.04827281913040844"
$ws.Range("K4").Value = [double]"3.0"
$ws.Range("L4").Value = [double]"1.0"
$ws.Range("M4").Value = [double]"0.478937"
$ws.Range("N4").Value = [double]"1.436811"
$ws.Range("O4").Value = [double]"0.0813355572127976"
$ws.Range("P4").Value = [double]"0.08133555721279762"
$ws.Range("Q4").Value = [double]"92.92579484862134"
$ws.Range("R4").Value = [double]"836.3321536375921"
$ws.Range("S4").Value = [double]"0.003926296642204365"
$ws.Range("T4").Value = [double]"0.003926296642204366"

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Bgn"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = [double]"3.0"
$ws.Range("F5").Value = [double]"1.0"
$ws.Range("G5").Value = [double]"194.0250906666667"
$ws.Range("H5").Value = [double]"582.075272"
$ws.Range("I5").Value = [double]"0.04827281913040843"
$ws.Range("J5").Value = [double]"0.04827281913040844"
$ws.Range("K5").Value = [double]"1.0"
$ws.Range("L5").Value = [double]"0.3333333333333333"
$ws.Range("M5").Value = [double]"0.1116203333333333"
$ws.Range("N5").Value = [double]"0.334861"
$ws.Range("O5").Value = [double]"0.01895594202983873"
$ws.Range("P5").Value = [double]"0.01895594202983874"
$ws.Range("Q5").Value = [double]"21.65714529524356"
$ws.Range("R5").Value = [double]"194.914307657192"
$ws.Range("S5").Value = [double]"0.0009150567610529123"
$ws.Range("T5").Value = [double]"0.0009150567610529127"

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bgn"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = [double]"3.0"
$ws.Range("F6").Value = [double]"1.0"
$ws.Range("G6").Value = [double]"3721.989909"
$ws.Range("H6").Value = [double]"11165.969727"
$ws.Range("I6").Value = [double]"0.9260191301290788"
$ws.Range("J6").Value = [double]"0.9260191301290789"
$ws.Range("K6").Value = [double]"3.0"
$ws.Range("L6").Value = [double]"1.0"
$ws.Range("M6").Value = [double]"4.618552666666667"
$ws.Range("N6").Value = [double]"13.855658"
$ws.Range("O6").Value = [double]"0.78434648953826"
$ws.Range("P6").Value = [double]"0.78434648953826"
$ws.Range("Q6").Value = [double]"17190.20641951837"
$ws.Range("R6").Value = [double]"154711.8577756654"
$ws.Range("S6").Value = [double]"0.7263198539620161"
$ws.Range("T6").Value = [double]"0.7263198539620163"

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bgn"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = [double]"3.0"
$ws.Range("F7").Value = [double]"1.0"
$ws.Range("G7").Value = [double]"3721.989909"
$ws.Range("H7").Value = [double]"11165.969727"
$ws.Range("I7").Value = [double]"0.9260191301290788"
$ws.Range("J7").Value = [double]"0.9260191301290789"
$ws.Range("K7").Value = [double]"2.0"
$ws.Range("L7").Value = [double]"0.6666666666666666"
$ws.Range("M7").Value = [double]"0.6792986666666666"
$ws.Range("N7").Value = [double]"2.037896"
$ws.Range("O7").Value = [double]"0.1153620112191035"
$ws.Range("P7").Value = [double]"0.1153620112191036"
$ws.Range("Q7").Value = [double]"2528.342782530488"
$ws.Range("R7").Value = [double]"22755.08504277439"
$ws.Range("S7").Value = [double]"0.1068274292790553"
$ws.Range("T7").Value = [double]"0.1068274292790553"

# Row 8: FAPs -> MuSCs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Bgn"
$ws.Range("C8").Value = "Fgfr3"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = [double]"3.0"
$ws.Range("F8").Value = [double]"1.0"
$ws.Range("G8").Value = [double]"3721.989909"
$ws.Range("H8").Value = [double]"11165.969727"
$ws.Range("I8").Value = [double]"0.9260191301290788"
$ws.Range("J8").Value = [double]"0.9260191301290789"
$ws.Range("K8").Value = [double]"3.0"
$ws.Range("L8").Value = [double]"1.0"
$ws.Range("M8").Value = [double]"0.478937"
$ws.Range("N8").Value = [double]"1.436811"
$ws.Range("O8").Value = [double]"0.0813355572127976"
$ws.Range("P8").Value = [double]"0.08133555721279762"
$ws.Range("Q8").Value = [double]"1782.598681046733"
$ws.Range("R8").Value = [double]"16043.3881294206"
$ws.Range("S8").Value = [double]"0.07531828193875875"
$ws.Range("T8").Value = [double]"0.07531828193875878"

# Row 9: FAPs -> Resolving-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Bgn"
$ws.Range("C9").Value = "Fgfr3"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = [double]"3.0"
$ws.Range("F9").Value = [double]"1.0"
$ws.Range("G9").Value = [double]"3721.989909"
$ws.Range("H9").Value = [double]"11165.969727"
$ws.Range("I9").Value = [double]"0.9260191301290788"
$ws.Range("J9").Value = [double]"0.9260191301290789"
$ws.Range("K9").Value = [double]"1.0"
$ws.Range("L9").Value = [double]"0.3333333333333333"
$ws.Range("M9").Value = [double]"0.1116203333333333"
$ws.Range("N9").Value = [double]"0.334861"
$ws.Range("O9").Value = [double]"0.01895594202983873"
$ws.Range("P9").Value = [double]"0.01895594202983874"
$ws.Range("Q9").Value = [double]"415.449754305883"
$ws.Range("R9").Value = [double]"3739.047788752947"
$ws.Range("S9").Value = [double]"0.01755356494924851"
$ws.Range("T9").Value = [double]"0.01755356494924851"

# Row 10: Inflammatory-Mac -> ECs
$ws.Range("A10").Value = "Inflammatory-Mac"
$ws.Range("B10").Value = "Bgn"
$ws.Range("C10").Value = "Fgfr3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = [double]"3.0"
$ws.Range("F10").Value = [double]"1.0"
$ws.Range("G10").Value = [double]"4.068420333333333"
$ws.Range("H10").Value = [double]"12.205261"
$ws.Range("I10").Value = [double]"0.001012209906577904"
$ws.Range("J10").Value = [double]"0.001012209906577904"
$ws.Range("K10").Value = [double]"3.0"
$ws.Range("L10").Value = [double]"1.0"
$ws.Range("M10").Value = [double]"4.618552666666667"
$ws.Range("N10").Value = [double]"13.855658"
$ws.Range("O10").Value = [double]"0.78434648953826"
$ws.Range("P10").Value = [double]"0.78434648953826"
$ws.Range("Q10").Value = [double]"18.79021357963756"
$ws.Range("R10").Value = [double]"169.111922216738"
$ws.Range("S10").Value = [double]"0.0007939232869002289"
$ws.Range("T10").Value = [double]"0.0007939232869002289"

# Row 11: Inflammatory-Mac -> FAPs
$ws.Range("A11").Value = "Inflammatory-Mac"
$ws.Range("B11").Value = "Bgn"
$ws.Range("C11").Value = "Fgfr3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = [double]"3.0"
$ws.Range("F11").Value = [double]"1.0"
$ws.Range("G11").Value = [double]"4.068420333333333"
$ws.Range("H11").Value = [double]"12.205261"
$ws.Range("I11").Value = [double]"0.001012209906577904"
$ws.Range("J11").Value = [double]"0.001012209906577904"
$ws.Range("K11").Value = [double]"2.0"
$ws.Range("L11").Value = [double]"0.6666666666666666"
$ws.Range("M11").Value = [double]"0.6792986666666666"
$ws.Range("N11").Value = [double]"2.037896"
$ws.Range("O11").Value = [double]"0.1153620112191035"
$ws.Range("P11").Value = [double]"0.1153620112191036"
$ws.Range("Q11").Value = [double]"2.763672507872888"
$ws.Range("R11").Value = [double]"24.873052570856"
$ws.Range("S11").Value = [double]"0.0001167705705987279"
$ws.Range("T11").Value = [double]"0.0001167705705987279"

# Row 12: Inflammatory-Mac -> MuSCs
$ws.Range("A12").Value = "Inflammatory-Mac"
$ws.Range("B12").Value = "Bgn"
$ws.Range("C12").Value = "Fgfr3"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = [double]"3.0"
$ws.Range("F12").Value = [double]"1.0"
$ws.Range("G12").Value = [double]"4.068420333333333"
$ws.Range("H12").Value = [double]"12.205261"
$ws.Range("I12").Value = [double]"0.001012209906577904"
$ws.Range("J12").Value = [double]"0.001012209906577904"
$ws.Range("K12").Value = [double]"3.0"
$ws.Range("L12").Value = [double]"1.0"
$ws.Range("M12").Value = [double]"0.478937"
$ws.Range("N12").Value = [double]"1.436811"
$ws.Range("O12").Value = [double]"0.0813355572127976"
$ws.Range("P12").Value = [double]"0.08133555721279762"
$ws.Range("Q12").Value = [double]"1.948517029185667"
$ws.Range("R12").Value = [double]"17.536653262671"
$ws.Range("S12").Value = [double]"8.23286567678276e-05"
$ws.Range("T12").Value = [double]"8.232865676782761e-05"

# Row 13: Inflammatory-Mac -> Resolving-Mac
$ws.Range("A13").Value = "Inflammatory-Mac"
$ws.Range("B13").Value = "Bgn"
$ws.Range("C13").Value = "Fgfr3"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = [double]"3.0"
$ws.Range("F13").Value = [double]"1.0"
$ws.Range("G13").Value = [double]"4.068420333333333"
$ws.Range("H13").Value = [double]"12.205261"
$ws.Range("I13").Value = [double]"0.001012209906577904"
$ws.Range("J13").Value = [double]"0.001012209906577904"
$ws.Range("K13").Value = [double]"1.0"
$ws.Range("L13").Value = [double]"0.3333333333333333"
$ws.Range("M13").Value = [double]"0.1116203333333333"
$ws.Range("N13").Value = [double]"0.334861"
$ws.Range("O13").Value = [double]"0.01895594202983873"
$ws.Range("P13").Value = [double]"0.01895594202983874"
$ws.Range("Q13").Value = [double]"0.4541184337467777"
$ws.Range("R13").Value = [double]"4.087065903721"
$ws.Range("S13").Value = [double]"1.918739231111922e-05"
$ws.Range("T13").Value = [double]"1.918739231111922e-05"

# Row 14: MuSCs -> ECs
$ws.Range("A14").Value = "MuSCs"
$ws.Range("B14").Value = "Bgn"
$ws.Range("C14").Value = "Fgfr3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = [double]"3.0"
$ws.Range("F14").Value = [double]"1.0"
$ws.Range("G14").Value = [double]"95.39550266666667"
$ws.Range("H14").Value = [double]"286.186508"
$ws.Range("I14").Value = [double]"0.02373409454550267"
$ws.Range("J14").Value = [double]"0.02373409454550267"
$ws.Range("K14").Value = [double]"3.0"
$ws.Range("L14").Value = [double]"1.0"
$ws.Range("M14").Value = [double]"4.618552666666667"
$ws.Range("N14").Value = [double]"13.855658"
$ws.Range("O14").Value = [double]"0.78434648953826"
$ws.Range("P14").Value = [double]"0.78434648953826"
$ws.Range("Q14").Value = [double]"440.5891532291405"
$ws.Range("R14").Value = [double]"3965.302379062264"
$ws.Range("S14").Value = [double]"0.01861575373913419"
$ws.Range("T14").Value = [double]"0.01861575373913419"

# Row 15: MuSCs -> FAPs
$ws.Range("A15").Value = "MuSCs"
$ws.Range("B15").Value = "Bgn"
$ws.Range("C15").Value = "Fgfr3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = [double]"3.0"
$ws.Range("F15").Value = [double]"1.0"
$ws.Range("G15").Value = [double]"95.39550266666667"
$ws.Range("H15").Value = [double]"286.186508"
$ws.Range("I15").Value = [double]"0.02373409454550267"
$ws.Range("J15").Value = [double]"0.02373409454550267"
$ws.Range("K15").Value = [double]"2.0"
$ws.Range("L15").Value = [double]"0.6666666666666666"
$ws.Range("M15").Value = [double]"0.6792986666666666"
$ws.Range("N15").Value = [double]"2.037896"
$ws.Range("O15").Value = [double]"0.1153620112191035"
$ws.Range("P15").Value = [double]"0.1153620112191036"
$ws.Range("Q15").Value = [double]"64.8020377674631"
$ws.Range("R15").Value = [double]"583.218339907168"
$ws.Range("S15").Value = [double]"0.002738012881233543"
$ws.Range("T15").Value = [double]"0.002738012881233544"

# Row 16: MuSCs -> MuSCs
$ws.Range("A16").Value = "MuSCs"
$ws.Range("B16").Value = "Bgn"
$ws.Range("C16").Value = "Fgfr3"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = [double]"3.0"
$ws.Range("F16").Value = [double]"1.0"
$ws.Range("G16").Value = [double]"95.39550266666667"
$ws.Range("H16").Value = [double]"286.186508"
$ws.Range("I16").Value = [double]"0.02373409454550267"
$ws.Range("J16").Value = [double]"0.02373409454550267"
$ws.Range("K16").Value = [double]"3.0"
$ws.Range("L16").Value = [double]"1.0"
$ws.Range("M16").Value = [double]"0.478937"
$ws.Range("N16").Value = [double]"1.436811"
$ws.Range("O16").Value = [double]"0.0813355572127976"
$ws.Range("P16").Value = [double]"0.08133555721279762"
$ws.Range("Q16").Value = [double]"45.68843586066534"
$ws.Range("R16").Value = [double]"411.195922745988"
$ws.Range("S16").Value = [double]"0.00193042580479968"
$ws.Range("T16").Value = [double]"0.001930425804799681"

# Row 17: MuSCs -> Resolving-Mac
$ws.Range("A17").Value = "MuSCs"
$ws.Range("B17").Value = "Bgn"
$ws.Range("C17").Value = "Fgfr3"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = [double]"3.0"
$ws.Range("F17").Value = [double]"1.0"
$ws.Range("G17").Value = [double]"95.39550266666667"
$ws.Range("H17").Value = [double]"286.186508"
$ws.Range("I17").Value = [double]"0.02373409454550267"
$ws.Range("J17").Value = [double]"0.02373409454550267"
$ws.Range("K17").Value = [double]"1.0"
$ws.Range("L17").Value = [double]"0.3333333333333333"
$ws.Range("M17").Value = [double]"0.1116203333333333"
$ws.Range("N17").Value = [double]"0.334861"
$ws.Range("O17").Value = [double]"0.01895594202983873"
$ws.Range("P17").Value = [double]"0.01895594202983874"
$ws.Range("Q17").Value = [double]"10.64807780615422"
$ws.Range("R17").Value = [double]"95.832700255388"
$ws.Range("S17").Value = [double]"0.0004499021203352603"
$ws.Range("T17").Value = [double]"0.0004499021203352604"

# Row 18: Resolving-Mac -> ECs
$ws.Range("A18").Value = "Resolving-Mac"
$ws.Range("B18").Value = "Bgn"
$ws.Range("C18").Value = "Fgfr3"
$ws.Range("D18").Value = "ECs"
$ws.Range("E18").Value = [double]"3.0"
$ws.Range("F18").Value = [double]"1.0"
$ws.Range("G18").Value = [double]"3.865589666666666"
$ws.Range("H18").Value = [double]"11.596769"
$ws.Range("I18").Value = [double]"0.0009617462884321383"
$ws.Range("J18").Value = [double]"0.0009617462884321383"
$ws.Range("K18").Value = [double]"3.0"
$ws.Range("L18").Value = [double]"1.0"
$ws.Range("M18").Value = [double]"4.618552666666667"
$ws.Range("N18").Value = [double]"13.855658"
$ws.Range("O18").Value = [double]"0.78434648953826"
$ws.Range("P18").Value = [double]"0.78434648953826"
$ws.Range("Q18").Value = [double]"17.85342946322244"
$ws.Range("R18").Value = [double]"160.680865169002"
$ws.Range("S18").Value = [double]"0.0007543423251581985"
$ws.Range("T18").Value = [double]"0.0007543423251581985"

# Row 19: Resolving-Mac -> FAPs
$ws.Range("A19").Value = "Resolving-Mac"
$ws.Range("B19").Value = "Bgn"
$ws.Range("C19").Value = "Fgfr3"
$ws.Range("D19").Value = "FAPs"
$ws.Range("E19").Value = [double]"3.0"
$ws.Range("F19").Value = [double]"1.0"
$ws.Range("G19").Value = [double]"3.865589666666666"
$ws.Range("H19").Value = [double]"11.596769"
$ws.Range("I19").Value = [double]"0.0009617462884321383"
$ws.Range("J19").Value = [double]"0.0009617462884321383"
$ws.Range("K19").Value = [double]"2.0"
$ws.Range("L19").Value = [double]"0.6666666666666666"
$ws.Range("M19").Value = [double]"0.6792986666666666"
$ws.Range("N19").Value = [double]"2.037896"
$ws.Range("O19").Value = [double]"0.1153620112191035"
$ws.Range("P19").Value = [double]"0.1153620112191036"
$ws.Range("Q19").Value = [double]"2.625889906447111"
$ws.Range("R19").Value = [double]"23.633009158024"
$ws.Range("S19").Value = [double]"0.0001109489861160395"
$ws.Range("T19").Value = [double]"0.0001109489861160395"

# Row 20: Resolving-Mac -> MuSCs
$ws.Range("A20").Value = "Resolving-Mac"
$ws.Range("B20").Value = "Bgn"
$ws.Range("C20").Value = "Fgfr3"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = [double]"3.0"
$ws.Range("F20").Value = [double]"1.0"
$ws.Range("G20").Value = [double]"3.865589666666666"
$ws.Range("H20").Value = [double]"11.596769"
$ws.Range("I20").Value = [double]"0.0009617462884321383"
$ws.Range("J20").Value = [double]"0.0009617462884321383"
$ws.Range("K20").Value = [double]"3.0"
$ws.Range("L20").Value = [double]"1.0"
$ws.Range("M20").Value = [double]"0.478937"
$ws.Range("N20").Value = [double]"1.436811"
$ws.Range("O20").Value = [double]"0.0813355572127976"
$ws.Range("P20").Value = [double]"0.08133555721279762"
$ws.Range("Q20").Value = [double]"1.851373918184333"
$ws.Range("R20").Value = [double]"16.662365263659"
$ws.Range("S20").Value = [double]"7.822417026696793e-05"
$ws.Range("T20").Value = [double]"7.822417026696794e-05"

# Row 21: Resolving-Mac -> Resolving-Mac
$ws.Range("A21").Value = "Resolving-Mac"
$ws.Range("B21").Value = "Bgn"
$ws.Range("C21").Value = "Fgfr3"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = [double]"3.0"
$ws.Range("F21").Value = [double]"1.0"
$ws.Range("G21").Value = [double]"3.865589666666666"
$ws.Range("H21").Value = [double]"11.596769"
$ws.Range("I21").Value = [double]"0.0009617462884321383"
$ws.Range("J21").Value = [double]"0.0009617462884321383"
$ws.Range("K21").Value = [double]"1.0"
$ws.Range("L21").Value = [double]"0.3333333333333333"
$ws.Range("M21").Value = [double]"0.1116203333333333"
$ws.Range("N21").Value = [double]"0.334861"
$ws.Range("O21").Value = [double]"0.01895594202983873"
$ws.Range("P21").Value = [double]"0.01895594202983874"
$ws.Range("Q21").Value = [double]"0.4314784071232222"
$ws.Range("R21").Value = [double]"3.883305664109"
$ws.Range("S21").Value = [double]"1.823080689093217e-05"
$ws.Range("T21").Value = [double]"1.823080689093218e-05"
